$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.212.99"
$ws.Range("E2").Value = "  -2.27%  "

$ws.Range("D3").Value = "1.872.75"
$ws.Range("E3").Value = "  -1.79%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5106"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3763"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07168"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8902"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07593"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.04%  "

$ws.Range("D13").Value = "1.858.50"
$ws.Range("E13").Value = "  -2.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.345"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008563"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.67%  "

$ws.Range("E18").Value = "  -2.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("D20").Value = "27.263.73"
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.081"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.71%  "

$ws.Range("D22").Value = "2.093.02"
$ws.Range("E22").Value = "  -1.80%  "

$ws.Range("E23").Value = "  -1.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.500"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.845"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.131"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.756"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.713"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08997"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05162"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.110"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7558"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.171"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02044"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.535"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.038"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("E40").Value = "  -1.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5363"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.664"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.581"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("E45").Value = "  -1.69%  "

$ws.Range("E46").Value = "  -2.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.38%  "

$ws.Range("E49").Value = "  -3.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "65.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
